# This workbook's data rows (2 and 4-12) got reshuffled: the "variable"
# columns (D, L, M, N, O, P, Q, R, S, T) of each row now hold the values
# that used to belong to a different row, per the mapping below
# (new row -> old row it now matches). Row 3, 13, 14 are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("D", "L", "M", "N", "O", "P", "Q", "R", "S", "T")
$rows = @(2, 4, 5, 6, 7, 8, 9, 10, 11, 12)

# Snapshot the current ("before") values for the columns that move.
# NOTE: use Value2 (not Value) -- Value does not resolve to the underlying
# primitive in this environment.
$snapshot = @{}
foreach ($r in $rows) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowVals
}

# new row -> source row (the row whose old values it should now contain)
$mapping = @{
    2  = 4
    4  = 5
    5  = 8
    6  = 9
    7  = 2
    8  = 11
    9  = 12
    10 = 6
    11 = 7
    12 = 10
}

foreach ($newRow in $mapping.Keys) {
    $srcRow = $mapping[$newRow]
    $srcVals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$newRow").Value2 = $srcVals[$c]
    }
}
